# Update posts.xlsx after post
# The row for the post "「休みが終わった」" (row 738) was removed from the
# spreadsheet; every row below it shifts up by one and the sheet's used
# range shrinks from A1:C819 to A1:C818.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Deleting the entire row shifts all subsequent rows up by one, which is
# exactly the transformation described by the diff (row 738 disappears,
# former row 739 becomes the new row 738, ..., former row 819 becomes the
# new row 818).
$ws.Rows(738).Delete()
